$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: insert a new paragraph "- Nhân viên" right before "- Người mua"
#         (paragraph 6 in the original document).
# ---------------------------------------------------------------------------
$pNguoiMua = $d.Paragraphs.Item(6)
$pNguoiMua.Range.InsertParagraphBefore()
$pNhanVien = $d.Paragraphs.Item(6)
$pNhanVien.Range.Text = "- Nhân viên"

# ---------------------------------------------------------------------------
# Step 2: delete the "Bảo mật:" paragraph entirely (now paragraph 8).
# ---------------------------------------------------------------------------
$pBaoMat = $d.Paragraphs.Item(8)
$pBaoMat.Range.Delete()

# ---------------------------------------------------------------------------
# Step 3: the former "- Đăng nhập có phân quyền" Heading1 (now paragraph 8)
#         becomes "Yêu cầu:" - keep the Heading1 style, just change the text.
# ---------------------------------------------------------------------------
$pHeading = $d.Paragraphs.Item(8)
$pHeading.Range.Text = "Yêu cầu:"

# ---------------------------------------------------------------------------
# Step 4: the former multi-run paragraph "- Người mua chỉ có thể đăng ký, ..."
#         (now paragraph 9) becomes a single run with new text.
# ---------------------------------------------------------------------------
$p9 = $d.Paragraphs.Item(9)
$r9 = $p9.Range
$r9.End = $r9.End - 1
$r9.Text = "- Ứng dụng sẽ hiển thị các sản phẩm, mô tả và đánh giá của sản phẩm cho người dùng xem, người dùng có thể chọn và thêm vào giỏ hàng, sau đó đặt mua, đặt mua xong mới có thể bình luận."

# ---------------------------------------------------------------------------
# Step 5: insert two brand-new paragraphs after paragraph 9.
#         The first one ("Cho phép nhiều phương thức...") needs two runs,
#         so it is built as two temporary paragraphs that get merged back
#         into one by deleting the paragraph mark between them.
# ---------------------------------------------------------------------------
$p9 = $d.Paragraphs.Item(9)
$p9.Range.InsertParagraphAfter()
$pPay1 = $d.Paragraphs.Item(10)
$rPay1 = $pPay1.Range
$rPay1.End = $rPay1.End - 1
$rPay1.Text = "- Cho phép nhiều phương thức thanh toán qua ví điện tử, tài khoản ngân hàng,"

$pPay1 = $d.Paragraphs.Item(10)
$pPay1.Range.InsertParagraphAfter()
$pPay2 = $d.Paragraphs.Item(11)
$rPay2 = $pPay2.Range
$rPay2.End = $rPay2.End - 1
$rPay2.Text = " ship COD..."

# merge paragraph 10 and 11 into a single paragraph with two runs
$pPay1 = $d.Paragraphs.Item(10)
$markStart = $pPay1.Range.End - 1
$markRange = $d.Range($markStart, $pPay1.Range.End)
$markRange.Delete()

# now insert the "- Mỗi người dùng cần đăng nhập..." paragraph after it
$pPay = $d.Paragraphs.Item(10)
$pPay.Range.InsertParagraphAfter()
$pLogin = $d.Paragraphs.Item(11)
$rLogin = $pLogin.Range
$rLogin.End = $rLogin.End - 1
$rLogin.Text = "- Mỗi người dùng cần đăng nhập, không có tài khoản thì sẽ đăng ký."

# ---------------------------------------------------------------------------
# Step 6: the former multi-run paragraph "- Admin được thêm, xóa, sửa ..."
#         (now paragraph 12) becomes a single run with new text.
# ---------------------------------------------------------------------------
$p12 = $d.Paragraphs.Item(12)
$r12 = $p12.Range
$r12.End = $r12.End - 1
$r12.Text = "- Nhân viên được phép xem các đơn chưa xử lý xong, đánh dấu thành đang ship hoặc đã ship"

# ---------------------------------------------------------------------------
# Step 7: the former "Phần cứng:" Heading1 (now paragraph 13) turns into a
#         plain (non-heading) paragraph with new text. Delete it and insert
#         a fresh plain paragraph in its place so no heading style lingers.
# ---------------------------------------------------------------------------
$p13 = $d.Paragraphs.Item(13)
$p13.Range.Delete()
$p12 = $d.Paragraphs.Item(12)
$p12.Range.InsertParagraphAfter()
$pShip = $d.Paragraphs.Item(13)
$rShip = $pShip.Range
$rShip.End = $rShip.End - 1
$rShip.Text = "- Các đơn đã ship sẽ được đưa vào lưu trữ."

# ---------------------------------------------------------------------------
# Step 8: "- JDK 1.8 trở lên" (now paragraph 14) becomes the Admin-rights text.
# ---------------------------------------------------------------------------
$p14 = $d.Paragraphs.Item(14)
$r14 = $p14.Range
$r14.End = $r14.End - 1
$r14.Text = "- Admin có tất cả quyền của nhân viên, nhưng có thể xóa đơn chưa xử lý xong, quản lý thông tin người dùng, quản lý các danh mục sản phẩm và các sản phẩm, xem thống kê."

# ---------------------------------------------------------------------------
# Step 9: "- Sử dụng JDBC" (now paragraph 15) becomes the statistics text.
# ---------------------------------------------------------------------------
$p15 = $d.Paragraphs.Item(15)
$r15 = $p15.Range
$r15.End = $r15.End - 1
$r15.Text = "- Thống kê gồm doanh thu theo tháng, theo năm, lượng đơn trên một người dùng."

# ---------------------------------------------------------------------------
# Step 10: delete the trailing "- SQL Server 2008 trở lên" paragraph (16).
# ---------------------------------------------------------------------------
$p16 = $d.Paragraphs.Item(16)
$p16.Range.Delete()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
